# ---------------------------------------------------------------------------
# Adds the "2022-Q4" quarter to the workbook:
#   1. Inserts a new row 2 into the "总计" (summary) sheet with the 2022-Q4
#      totals, shifting the existing quarters down by one row.
#   2. Inserts a brand-new worksheet named "2022-Q4" (with the fund detail
#      table) right after "总计", shifting all the later quarter sheets down
#      by one tab position.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Helper: assign a value to a cell while forcing a "numeric-looking" string
# (e.g. "18.20") to be stored as TEXT, just like the rest of the workbook's
# fund-detail tables. Using a NumberFormat round-trip keeps the shared-string
# type without leaving a lingering cell style behind.
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1. "总计" summary sheet: insert the new 2022-Q4 row.
# ---------------------------------------------------------------------------

$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 9
$wsTotal.Range("D2").Value = 0.75

# Carry the bold/bordered "index column" style from the row below onto the
# freshly inserted A2 cell (Insert() does not copy it), and make sure
# B2:D2 keep the plain/default style (Insert() does not add one, but stay
# defensive in case the host copies formatting from the row above).
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$wsTotal.Range("B2:D2").ClearFormats()

# Re-apply cell values after the formatting operations above.
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 9
$wsTotal.Range("D2").Value = 0.75

Write-Host "总计 sheet updated"

# ---------------------------------------------------------------------------
# 2. Brand-new "2022-Q4" fund-detail worksheet.
# ---------------------------------------------------------------------------

$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q4"
$newSheet.Move($wsQ3)

# NOTE: the $newSheet handle can go stale (rebind to whatever sheet now
# occupies its old slot) once Move() reshuffles tab positions, so re-fetch
# a fresh reference by name before writing any data into it.
$newSheet = $wb.Worksheets.Item("2022-Q4")

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows: A (index, numeric), B (fund code, text), C (fund name, text),
# D/E/F/G (numeric-looking text values), H (rank, numeric).
$rows = @(
    @(0, "502000", "西部利得中证500指数增强（LOF）A", "18.20", "90.25", "1.50", "0.2730", 10),
    @(1, "005994", "国投瑞银中证500指数量化增强A",     "13.06", "89.93", "1.41", "0.1841", 7),
    @(2, "014155", "国泰君安中证500指数增强A",         "7.70",  "92.93", "1.06", "0.0816", 6),
    @(3, "009300", "西部利得中证500指数增强（LOF）C", "4.99",  "90.25", "1.50", "0.0748", 10),
    @(4, "007089", "国投瑞银中证500指数量化增强C",     "3.73",  "89.93", "1.41", "0.0526", 7),
    @(5, "014156", "国泰君安中证500指数增强C",         "4.81",  "92.93", "1.06", "0.0510", 6),
    @(6, "540004", "汇丰晋信2026周期混合",             "1.08",  "23.97", "1.10", "0.0119", 9),
    @(7, "008778", "嘉实中证500指数增强A",             "0.60",  "93.52", "1.81", "0.0109", 6),
    @(8, "008779", "嘉实中证500指数增强C",             "0.40",  "93.52", "1.81", "0.0072", 6)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $newSheet.Range("A$r").Value = $data[0]
    Set-TextValue $newSheet.Range("B$r") $data[1]
    Set-TextValue $newSheet.Range("C$r") $data[2]
    Set-TextValue $newSheet.Range("D$r") $data[3]
    Set-TextValue $newSheet.Range("E$r") $data[4]
    Set-TextValue $newSheet.Range("F$r") $data[5]
    Set-TextValue $newSheet.Range("G$r") $data[6]
    $newSheet.Range("H$r").Value = $data[7]
}

Write-Host "2022-Q4 sheet created"
Write-Host "Done"
